$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 (shifts existing rows 4-9 down to 5-10),
# copying formatting from the row that was previously row 4.
$ws.Rows(4).Insert()

# Populate the new row 4 with the new SNP-count dataset.
$ws.Range("A4").Value = "Corallorhiza bentleyi -ginger dev"
$ws.Range("B4").Value = "ISSRseq"
$ws.Range("C4").Value = "pooled PCRs, sheared"
$ws.Range("D4").Value = "de novo"
$ws.Range("E4").Value = "3,185 (250 bp min)"
$ws.Range("F4").Value = 1435470
$ws.Range("G4").Value = 130540
$ws.Range("H4").Value = 40447
$ws.Range("I4").Value = 32712
$ws.Range("J4").Value = 11942

# Highlight the new row with a yellow fill, like the rest of the row's
# formatting, to flag it as the newest development run.
$ws.Range("A4:J4").Interior.Color = 65535

# Match the saved selection/view state from the edited workbook.
$ws.Range("J4").Select()
